$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Collapse the run split around "...into even more detail. That is, "
#    and "break down the output..." into a single run (same visible text,
#    just merges two adjacent <w:r> elements into one).
# ---------------------------------------------------------------------
$mergeRange = $d.Content
$mergeRange.Find.Text = "break down the output not only into parts of speech (which is determined by the first letter of the code/tag), but also into the tenses/moods/person etc, (represented by each individual letter of the entire tag)."
$mergeRange.Find.Execute() | Out-Null
$afterPhraseStart = $mergeRange.Start

$beforeRange = $d.Content
$beforeRange.Find.Text = " into even more detail. That is, "
$beforeRange.Find.Execute() | Out-Null
$beforePhraseStart = $beforeRange.Start

$fullRange = $d.Range($beforePhraseStart, $afterPhraseStart + 255)
# Recompute the precise end using the matched range from the second find
$fullRange = $d.Range($beforePhraseStart, $mergeRange.End)

$mergedText = " into even more detail. That is, break down the output not only into parts of speech (which is determined by the first letter of the code/tag), but also into the tenses/moods/person etc, (represented by each individual letter of the entire tag)."
$fullRange.Text = ""
$insertPoint = $d.Range($beforePhraseStart, $beforePhraseStart)
$insertPoint.InsertAfter($mergedText)

# ---------------------------------------------------------------------
# 2. Append a new paragraph at the end of the document with the new
#    narrative text about using Rapid Miner / Excel, preserved as three
#    separate runs (matching the source run-split).
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$tailRange = $lastPara.Range
$tailRange.Collapse(0)
$tailRange.InsertAfter([char]13)

$newPara = $d.Paragraphs.Last
$newRange = $newPara.Range
$newRange.Collapse(1)

$rsquo = [char]0x2019
$newParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">For instance, I used Rapid Miner to create a breakdown of </w:t></w:r><w:r><w:t>the frequency of verbs by mood, tense, person and number</w:t></w:r><w:r><w:t xml:space="preserve"> and put the results in Analysis.xlsx. I decided to use Excel since that is a simple way for anyone to do basic statistical analysis. One can sort by the field in which they' + $rsquo + 're interested, or sum the percentages to see what percent of the entire corpus is comprised of whichever subset they want.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newRange.InsertXML($newParaXml)

# InsertXML of a full <w:p> leaves behind the now-empty paragraph mark
# that used to terminate the document; remove that extra trailing
# paragraph break so the inserted paragraph becomes the document's last
# paragraph again (mirrors the original body/sectPr structure).
$trailingPara = $d.Paragraphs.Last
$delStart = $trailingPara.Range.Start - 1
$delEnd = $trailingPara.Range.End
$d.Range($delStart, $delEnd).Delete()

# ---------------------------------------------------------------------
# 3. Relocate the "_GoBack" bookmark so it again marks the end of the
#    document's (now new) last paragraph, i.e. right after the text and
#    before the paragraph mark -- same relative position it had before
#    the new paragraph was appended.
# ---------------------------------------------------------------------
$finalPara = $d.Paragraphs.Last
$bmPos = $finalPara.Range.End - 1

# A collapsed Range sitting exactly on "end-of-text, before the pilcrow"
# trips a positioning quirk when handed straight to Bookmarks.Add, so we
# nudge it: insert a throwaway character at that spot, anchor the
# bookmark before the throwaway character, then remove the character.
$nudge = $d.Range($bmPos, $bmPos)
$nudge.InsertAfter("Z")
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
$d.Range($bmPos, $bmPos + 1).Delete()
